# Update SCD0026-001 until SCD0026-017 Fix
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the sheet tab from SCD0338 to SCD0026
$ws.Name = "SCD0026"

# Update the TC_ID cell from SCD0338-015 to SCD0026-015
$ws.Range("B2").Value = "SCD0026-015"

# Reset view: scroll back to A1 and move selection to B3
[void]$ws.Range("A1").Select()
[void]$ws.Range("B3").Select()
